$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.913.85'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.633.43'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.82'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5116'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2565'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06347'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.50'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07785'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.252'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '1.631.34'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '1.858.64'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5523'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.82'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '0.0₅7610'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '25.904.94'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.23'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.866'
$ws.Range('E22').Value = '  -0.57%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.890'
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.16'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1258'
$ws.Range('E27').Value = '  +5.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.770'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.242'
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04916'
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.238'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.187'
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8978'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5528'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.541'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').Value = '1.116.00'
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01555'
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.590'
$ws.Range('E42').Value = '  +3.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7949'
$ws.Range('E43').Value = '  -1.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.79'
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').Value = '1.768.94'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = '0.0₈117'
$ws.Range('E46').Value = '  -8.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4437'
$ws.Range('E47').Value = '  -2.05%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.86'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('E50').Value = '  +1.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.570'
$ws.Range('E51').Value = '  +3.25%  '
